# Update the heading date line.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2023-09-01 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-09-02 Saturday", 2)

# Update the division problems in the table, cell-by-cell so that
# duplicate expressions (e.g. "79÷4=" which occurs three times with
# three different replacements) are each handled independently based
# on their (row, column) position rather than a global text match.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "70÷4="
$t.Cell(1,2).Range.Text  = "82÷2="
$t.Cell(1,3).Range.Text  = "72÷2="
$t.Cell(1,4).Range.Text  = "79÷4="
$t.Cell(1,5).Range.Text  = "27÷6="

$t.Cell(5,1).Range.Text  = "81÷4="
$t.Cell(5,2).Range.Text  = "45÷4="
$t.Cell(5,3).Range.Text  = "64÷9="
$t.Cell(5,4).Range.Text  = "27÷4="
$t.Cell(5,5).Range.Text  = "43÷3="

$t.Cell(9,1).Range.Text  = "44÷7="
$t.Cell(9,2).Range.Text  = "68÷5="
$t.Cell(9,3).Range.Text  = "28÷5="
$t.Cell(9,4).Range.Text  = "79÷8="
$t.Cell(9,5).Range.Text  = "94÷3="

$t.Cell(13,1).Range.Text = "14÷7="
$t.Cell(13,2).Range.Text = "69÷2="
$t.Cell(13,3).Range.Text = "57÷4="
$t.Cell(13,4).Range.Text = "84÷3="
$t.Cell(13,5).Range.Text = "57÷2="

$t.Cell(17,1).Range.Text = "77÷8="
$t.Cell(17,2).Range.Text = "45÷5="
$t.Cell(17,3).Range.Text = "45÷6="
$t.Cell(17,4).Range.Text = "46÷5="
$t.Cell(17,5).Range.Text = "81÷4="
